$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "language" column (F) with a header and a value per data row,
# in top-to-bottom order so the new shared-string entries are created in
# the same order the author's edit produced them (language, Russian,
# Belarussian, Both, ...). Row 16 and row 29 are intentionally skipped -
# those rows have no language value, matching the source edit.
$languageByRow = @(
    ,@(1, "language")
    ,@(2, "Russian")
    ,@(3, "Russian")
    ,@(4, "Russian")
    ,@(5, "Russian")
    ,@(6, "Russian")
    ,@(7, "Belarussian")
    ,@(8, "Russian")
    ,@(9, "Belarussian")
    ,@(10, "Russian")
    ,@(11, "Russian")
    ,@(12, "Belarussian")
    ,@(13, "Belarussian")
    ,@(14, "Russian")
    ,@(15, "Russian")
    ,@(17, "Russian")
    ,@(18, "Both")
    ,@(19, "Both")
    ,@(20, "Both")
    ,@(21, "Russian")
    ,@(22, "Russian")
    ,@(23, "Russian")
    ,@(24, "Russian")
    ,@(25, "Russian")
    ,@(26, "Belarussian")
    ,@(27, "Russian")
    ,@(28, "Belarussian")
)

foreach ($pair in $languageByRow) {
    $row = $pair[0]
    $value = $pair[1]
    $ws.Cells.Item($row, 6).Value = $value
}

# Clean up the inconsistent handle text: "studenty_by" -> "studentyBY"
$ws.Range("B24").Value = "studentyBY"

# Match the author's final selection in the sheet view
$ws.Range("B10").Select()
